$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete columns I (ExtractionType) and J (SamplePortion) entirely,
# shifting the old column K (Comment) left to become column I.
$ws.Range("I1:J3").Delete()
